# edit.ps1 - applies the tracked-change edits described by the diff to
# the geodetic-point reference table in the document.
#
# Word find/replace constants used below:
#   wdReplaceNone = 0, wdReplaceOne = 1, wdReplaceAll = 2
# We consistently use wdReplaceOne (1) so that $range.Start/.End collapse
# onto the single match just made (needed for the scoped "клacc" swaps).

$d = $word.ActiveDocument

function Replace-UniqueText($old, $new) {
    $r = $d.Content
    $r.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 1) | Out-Null
}

# ---------------------------------------------------------------------
# 1) The two ambiguous "N класс" cells must swap values with each other
#    (1509-row: 3->2, 1525-row: 2->3), while the 1467-row's own
#    "2 класс" must stay untouched. We scope each Find to a narrow
#    range bounded by the still-unique point numbers that precede each
#    row, so the correct occurrence is hit regardless of the other
#    edits happening elsewhere in the document.
# ---------------------------------------------------------------------

# 1509 Пролетарский row: "3 класс" -> "2 класс" (this text is globally
# unique before any edits, so no extra scoping is required).
$rClass1 = $d.Content
$rClass1.Find.Execute("3 класс", $false, $false, $false, $false, $false, $true, 1, $false, "2 класс", 1) | Out-Null

# 1525 Мазово row: "2 класс" -> "3 класс" (first of two remaining
# "2 класс" occurrences). Scope the search between the "1525" and
# "1467" anchors so the 1467-row's "2 класс" is left alone.
$bound1 = $d.Content
$bound1.Find.Execute("1525", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bound2 = $d.Content
$bound2.Find.Execute("1467", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$scopedClass = $d.Range($bound1.Start, $bound2.Start)
$scopedClass.Find.Execute("2 класс", $false, $false, $false, $false, $false, $true, 1, $false, "3 класс", 1) | Out-Null

# ---------------------------------------------------------------------
# 2) Row 2 ("№ п/п" = 2, empty document-name cell) gains a document
#    title. The cell currently holds two runs: a " " run followed by an
#    empty run. We rewrite the whole cell's visible text (minus the
#    trailing paragraph mark) so the final text reads
#    " Постановление Администрации Торжокского района Тверской области".
# ---------------------------------------------------------------------

$table = $d.Tables.Item(1)
$docNameCell = $table.Rows.Item(6).Cells.Item(2)
$cellRange = $docNameCell.Range
$cellRange.End = $cellRange.End - 1
$cellRange.Text = " Постановление Администрации Торжокского района Тверской области"

# ---------------------------------------------------------------------
# 3) Remaining cells each hold one unique run of text in the document,
#    so a direct global Find/Replace is safe for each of them.
# ---------------------------------------------------------------------

Replace-UniqueText " № 02-69/17-1-89716 от 17.02.2017 г." " № 99/2017/31874972 от 19.10.2017 г."
Replace-UniqueText " № 99/2017/37729316 от 22.11.2017 г." " № 96 от 16.02.2018 г."
Replace-UniqueText "Постановление Администрации Старицкого района Тверской области" "Выписка из каталога геодезических пунктов на Тверскую область"
Replace-UniqueText " № 179 от 06.04.2018 г." " № б/н от 28.02.2018 г."

Replace-UniqueText "1509 Пролетарский пир. 6 " "1313 Владенино пирамида "
Replace-UniqueText "375354.6" "319215.36"
Replace-UniqueText "2191497.81" "2222474.13"

Replace-UniqueText "1525 Мазово наружный знак утрачен " "1300  Колодези пирамида "
Replace-UniqueText "379294.24" "315734.2"
Replace-UniqueText "2223444.25" "2224862.19"

Replace-UniqueText "1467 Еваново наружный знак утрачен " "1285 Торжок наружный знак утрачен "
Replace-UniqueText "362953.1" "311784.21"
Replace-UniqueText "2196831.21" "2215334.98"

Replace-UniqueText "GNSS - приемник спутниковый геодезический двухчастотный Trimble R8" "GNSS - приёмник спутниковый геодезический двухчастотный Trimble R8 GNSS"
Replace-UniqueText "33967-07 до 10.01.2018г." "33967-07 до 17.01.2018г."
Replace-UniqueText "№ 012343    " "№012343 от 18.01.2017г."

Replace-UniqueText "Аппаратура геодезическая спутниковая Stonex S9 GNSS" "Аппаратура геодезическая потребителей спутниковых навигационных систем ГЛОНАСС и GPS Trimble R7GNSS"
Replace-UniqueText "№ 50874-12 до 13.03.2018г." "37145-08 до 17.01.2018г."
Replace-UniqueText "№ 013364   " "№012342 от 18.01.2017г."

Write-Output "edits applied"
